{"js": "// The document contains a paragraph whose text is the M2Doc user-doc\n// opening field marker: {m:userdoc 'zone1'}\n// It is currently split across 2 runs (\"{m\" and \":userdoc 'zone1'}\").\n// The commit migrates the parser to split field markers into 4 runs:\n// \"{\", \"m\", \":userdoc 'zone1'\", \"}\" (one run per logical token piece).\n//\n// We locate the exact text, then replace it in-place with an OOXML\n// fragment that contains the four separate <w:r> runs so the split\n// survives (a plain insertText() call would just re-merge adjacent\n// runs that share formatting).\n\nconst searchText = \"{m:userdoc 'zone1'}\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text: \" + searchText);\n}\n\nconst target = results.items[0];\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:userdoc \\'zone1\\'</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document contains a paragraph whose text is the M2Doc user-doc\n# opening field marker: {m:userdoc 'zone1'}\n# It is currently split across 2 runs (\"{m\" and \":userdoc 'zone1'}\").\n# The commit migrates the parser to split field markers into 4 runs:\n# \"{\", \"m\", \":userdoc 'zone1'\", \"}\" (one run per logical token piece).\n#\n# Word COM's Range.InsertXML on this host inserts the supplied OOXML\n# right after the target range (it does not replace the range's\n# contents), so we insert the 4-run replacement immediately before the\n# original text, then delete the now-shifted original text. This keeps\n# the host paragraph's own attributes untouched (only the runs inside\n# it change), matching a plain run re-split.\n\n$d = $word.ActiveDocument\n\n$searchText = \"{m:userdoc 'zone1'}\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute($searchText)\nif (-not $found) {\n    throw \"Could not find target text: $searchText\"\n}\n\n$start = $rng.Start\n$end = $rng.End\n$length = $end - $start\n\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">\n<pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:userdoc 'zone1'</w:t></w:r><w:r><w:t xml:space=\"preserve\">}</w:t></w:r></w:p></w:body></w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n# Insert the 4-run replacement right before the original text.\n$insertPoint = $d.Range($start, $start)\n$insertPoint.InsertXML($ooxml)\n\n# The original text got pushed later by the inserted content; remove it.\n$oldRange = $d.Range($start + $length, $end + $length)\n$oldRange.Text = \"\"\n"}
